$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3884.7273
$ws.Range("J17").Value = 2766.9
$ws.Range("L17").Value = 8300.700000000001
$ws.Range("N17").Value = -8636.700000000001

$ws.Range("H28").Value = 1272.9
$ws.Range("I28").Value = 995.3077
$ws.Range("K28").Value = 995.3077
$ws.Range("M28").Value = -510.3077

$ws.Range("H86").Value = 4148.8184
$ws.Range("I86").Value = 3981.6667
$ws.Range("K86").Value = 3981.6667
$ws.Range("M86").Value = -2858.6667

$ws.Range("I88").Value = 5599.375
$ws.Range("J88").Value = 2001.3334
$ws.Range("K88").Value = 5599.375
$ws.Range("L88").Value = 2001.3334
$ws.Range("M88").Value = -5193.375
$ws.Range("N88").Value = -2813.3334

$ws.Range("H89").Value = 4148.8184
$ws.Range("I89").Value = 3981.6667
$ws.Range("K89").Value = 19908.3335
$ws.Range("M89").Value = -14292.3335

$ws.Range("I91").Value = 5599.375
$ws.Range("J91").Value = 2001.3334
$ws.Range("K91").Value = 5599.375
$ws.Range("L91").Value = 2001.3334
$ws.Range("M91").Value = -4195.375
$ws.Range("N91").Value = -4809.3334

$ws.Range("H92").Value = 970.7727
$ws.Range("I92").Value = 767.7143
$ws.Range("K92").Value = 767.7143
$ws.Range("M92").Value = 480.2857

$ws.Range("H104").Value = 118.85714
$ws.Range("I104").Value = 118.85714
$ws.Range("K104").Value = 356.57142
$ws.Range("M104").Value = 1390.42858

$ws.Range("H113").Value = 8251
$ws.Range("I113").Value = 8321
$ws.Range("J113").Value = 8221.833000000001
$ws.Range("K113").Value = 8321
$ws.Range("L113").Value = 8221.833000000001
$ws.Range("M113").Value = -5067
$ws.Range("N113").Value = -14729.833

$ws.Range("H132").Value = 13416.308
$ws.Range("I132").Value = 3165.4285
$ws.Range("K132").Value = 9496.2855
$ws.Range("M132").Value = -6966.2855

$ws.Range("H135").Value = 11113951
$ws.Range("I135").Value = 12196386
$ws.Range("K135").Value = 109767474
$ws.Range("M135").Value = -109764939

$ws.Range("H138").Value = 3230.9656
$ws.Range("J138").Value = 4002.8572
$ws.Range("L138").Value = 12008.5716
$ws.Range("N138").Value = -22288.5716

$ws.Range("H141").Value = 3935.9092
$ws.Range("I141").Value = 2088.611
$ws.Range("J141").Value = 12248.75
$ws.Range("K141").Value = 6265.833
$ws.Range("L141").Value = 36746.25
$ws.Range("M141").Value = -1085.833
$ws.Range("N141").Value = -47106.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 5083.1665
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 4099.8
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 4099.8
$ws.Range("M50").Value = -9286
$ws.Range("N50").Value = -5527.8

$ws.Range("H74").Value = 2944.111
$ws.Range("I74").Value = 2374.625
$ws.Range("K74").Value = 2374.625
$ws.Range("M74").Value = -1500.625

$ws.Range("H77").Value = 2944.111
$ws.Range("I77").Value = 2374.625
$ws.Range("K77").Value = 11873.125
$ws.Range("M77").Value = -7505.125

$ws.Range("H110").Value = 907.8095
$ws.Range("I110").Value = 959.1111
$ws.Range("K110").Value = 959.1111
$ws.Range("M110").Value = 1085.8889

$ws.Range("H132").Value = 940.5517
$ws.Range("I132").Value = 680.4091
$ws.Range("J132").Value = 1758.1428
$ws.Range("K132").Value = 2041.2273
$ws.Range("L132").Value = 5274.428400000001
$ws.Range("M132").Value = 488.7727
$ws.Range("N132").Value = -10334.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 60966.824
$ws.Range("I20").Value = 1338.6
$ws.Range("K20").Value = 1338.6
$ws.Range("M20").Value = -1091.6

$ws.Range("H86").Value = 2667.1282
$ws.Range("J86").Value = 3082.1
$ws.Range("L86").Value = 3082.1
$ws.Range("N86").Value = -5328.1

$ws.Range("H89").Value = 2667.1282
$ws.Range("J89").Value = 3082.1
$ws.Range("L89").Value = 15410.5
$ws.Range("N89").Value = -26642.5

$ws.Range("H134").Value = 2077.8696
$ws.Range("I134").Value = 1585.881
$ws.Range("J134").Value = 7243.75
$ws.Range("K134").Value = 4757.643
$ws.Range("L134").Value = 21731.25
$ws.Range("M134").Value = -2222.643
$ws.Range("N134").Value = -26801.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2137.647
$ws.Range("I58").Value = 1622.7333
$ws.Range("J58").Value = 5999.5
$ws.Range("K58").Value = 1622.7333
$ws.Range("L58").Value = 5999.5
$ws.Range("M58").Value = -1419.7333
$ws.Range("N58").Value = -6405.5

$ws.Range("H93").Value = 26118
$ws.Range("I93").Value = 3789.8
$ws.Range("J93").Value = 63331.668
$ws.Range("K93").Value = 3789.8
$ws.Range("L93").Value = 63331.668
$ws.Range("M93").Value = -1917.8
$ws.Range("N93").Value = -67075.66800000001

$ws.Range("H99").Value = 8633600
$ws.Range("I99").Value = 1744323.9
$ws.Range("K99").Value = 1744323.9
$ws.Range("M99").Value = -1742825.9

$ws.Range("H100").Value = 80780
$ws.Range("J100").Value = 80780
$ws.Range("L100").Value = 80780
$ws.Range("N100").Value = -82944

$ws.Range("H126").Value = 8633600
$ws.Range("I126").Value = 1744323.9
$ws.Range("K126").Value = 5232971.699999999
$ws.Range("M126").Value = -5230501.699999999

$ws.Range("H132").Value = 2066.5833
$ws.Range("I132").Value = 1891.1212
$ws.Range("K132").Value = 5673.363600000001
$ws.Range("M132").Value = -3143.363600000001

$ws.Range("H134").Value = 3008.465
$ws.Range("I134").Value = 2686.5676
$ws.Range("J134").Value = 4993.5
$ws.Range("K134").Value = 8059.702799999999
$ws.Range("L134").Value = 14980.5
$ws.Range("M134").Value = -5524.702799999999
$ws.Range("N134").Value = -20050.5

$ws.Range("H136").Value = 2137.647
$ws.Range("I136").Value = 1622.7333
$ws.Range("J136").Value = 5999.5
$ws.Range("K136").Value = 4868.199900000001
$ws.Range("L136").Value = 17998.5
$ws.Range("M136").Value = -2318.199900000001
$ws.Range("N136").Value = -23098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 158.73914
$ws.Range("I61").Value = 123.35
$ws.Range("J61").Value = 394.66666
$ws.Range("K61").Value = 370.05
$ws.Range("L61").Value = 1183.99998
$ws.Range("M61").Value = -155.05
$ws.Range("N61").Value = -1613.99998

$ws.Range("H80").Value = 4476
$ws.Range("J80").Value = 4666.6
$ws.Range("L80").Value = 13999.8
$ws.Range("N80").Value = -15871.8

$ws.Range("H83").Value = 4476
$ws.Range("J83").Value = 4666.6
$ws.Range("L83").Value = 41999.4
$ws.Range("N83").Value = -51359.4

$ws.Range("H113").Value = 1751.5
$ws.Range("I113").Value = 2362
$ws.Range("K113").Value = 7086
$ws.Range("M113").Value = -4916

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 83629.21000000001
$ws.Range("I80").Value = 128539.555
$ws.Range("K80").Value = 128539.555
$ws.Range("M80").Value = -127541.555

$ws.Range("H83").Value = 83629.21000000001
$ws.Range("I83").Value = 128539.555
$ws.Range("K83").Value = 642697.7749999999
$ws.Range("M83").Value = -637705.7749999999

$ws.Range("H97").Value = 13842.143
$ws.Range("J97").Value = 34000
$ws.Range("L97").Value = 34000
$ws.Range("N97").Value = -34992

$ws.Range("H132").Value = 2122.2188
$ws.Range("I132").Value = 1858.5518
$ws.Range("K132").Value = 5575.6554
$ws.Range("M132").Value = -3045.6554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2085.5715
$ws.Range("I22").Value = 1933
$ws.Range("J22").Value = 2200
$ws.Range("K22").Value = 1933
$ws.Range("L22").Value = 2200
$ws.Range("M22").Value = -1638
$ws.Range("N22").Value = -2790

$ws.Range("H27").Value = 2085.5715
$ws.Range("I27").Value = 1933
$ws.Range("J27").Value = 2200
$ws.Range("K27").Value = 1933
$ws.Range("L27").Value = 2200
$ws.Range("M27").Value = -1826
$ws.Range("N27").Value = -2414

$ws.Range("H46").Value = 3944.4666
$ws.Range("J46").Value = 3931.8408
$ws.Range("L46").Value = 3931.8408
$ws.Range("N46").Value = -4307.8408

$ws.Range("H55").Value = 2438.2632
$ws.Range("I55").Value = 3443.2222
$ws.Range("K55").Value = 3443.2222
$ws.Range("M55").Value = -3270.2222

$ws.Range("H132").Value = 4962.514
$ws.Range("I132").Value = 3381.3809
$ws.Range("J132").Value = 7334.2144
$ws.Range("K132").Value = 10144.1427
$ws.Range("L132").Value = 22002.6432
$ws.Range("M132").Value = -7614.1427
$ws.Range("N132").Value = -27062.6432

$ws.Range("H136").Value = 4726.857
$ws.Range("I136").Value = 4579.684
$ws.Range("J136").Value = 6125
$ws.Range("K136").Value = 13739.052
$ws.Range("L136").Value = 18375
$ws.Range("M136").Value = -11189.052
$ws.Range("N136").Value = -23475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 41747.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 41747.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 41747.5
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -42207.5

$ws.Range("H76").Value = 42500
$ws.Range("I76").Value = 40000
$ws.Range("J76").Value = 45000
$ws.Range("K76").Value = 40000
$ws.Range("L76").Value = 45000
$ws.Range("M76").Value = -39685
$ws.Range("N76").Value = -45630

$ws.Range("H79").Value = 42500
$ws.Range("I79").Value = 40000
$ws.Range("J79").Value = 45000
$ws.Range("K79").Value = 40000
$ws.Range("L79").Value = 45000
$ws.Range("M79").Value = -38908
$ws.Range("N79").Value = -47184

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("M123").Value = -59800

$ws.Range("H126").Value = 2423
$ws.Range("I126").Value = 2100.875
$ws.Range("K126").Value = 6302.625
$ws.Range("M126").Value = -3832.625

$ws.Range("H132").Value = 1603.3684
$ws.Range("I132").Value = 1414.6666
$ws.Range("K132").Value = 4243.9998
$ws.Range("M132").Value = -1713.9998
